$wb = $excel.ActiveWorkbook

# --- Add a new "Player Info" sheet in front of the existing sheets ---
$ws = $wb.Worksheets.Add()
$ws.Name = "Player Info"

# Header row: bold, centered, top-aligned, thin box border (matches the
# header style already used on the other two sheets).
$hdr = $ws.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "NAME"
$ws.Range("C1").Value = "BATTING_HAND"
$ws.Range("D1").Value = "BOWL_STYLE"

# Data row. ID is a numeric-looking code, but must be stored as text, so
# force the text format before writing the value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4815"
$ws.Range("B2").Value = "Sherman Hakim Lewis"
$ws.Range("C2").Value = "Right Handed"
$ws.Range("D2").Value = "Right Arm Fast"

# --- Rename MATCH_CARD_LINK -> MATCH_CODE and store the bare match code
#     instead of the full scorecard URL, on both remaining sheets. ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4583"

$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4583"
